# Append 6 new weekly scoreboard rows (rows 318-323) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data as captured from the source diff: Participant, Date(serial),
# Workout Type, Total Duration, Total Distance, Total Elevation,
# Zone 1..5, Workout Level, Week
$rows = @(
    @("Matt",     45505, "Run",     33, 4,     256, 0,  2,  9,  9,  12, "Agile Antelope",   8),
    @("Steven",   45505, "Walk",    23, 1.09,  36,  23, 0,  0,  0,  0,  "Brave Leopard",    8),
    @("Phil",     45505, "Workout", 69, 0,     0,   38, 21, 6,  4,  0,  "Sauntering Hippo", 8),
    @("Jeremiah", 45506, "Workout", 68, 0,     0,   51, 17, 0,  0,  0,  "Agile Antelope",   8),
    @("Steven",   45506, "Run",     39, 3.26,  95,  18, 21, 1,  0,  0,  "Brave Leopard",    8),
    @("Eric",     45506, "Run",     77, 7.3,   223, 0,  31, 43, 0,  0,  "Agile Antelope",   8)
)

# Template cell carrying the date number format already used throughout
# column B (style index referencing numFmtId 14 = m/d/yyyy), so the copy
# below reuses the existing style instead of minting a new one.
$dateTemplate = $ws.Cells.Item(317, 2)

$startRow = 318
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]

    $dateCell = $ws.Cells.Item($r, 2)
    $dateTemplate.Copy($dateCell)
    $dateCell.Value = $data[1]

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
    $ws.Cells.Item($r, 13).Value = $data[12]
}

$ws.Range("A324").Select() | Out-Null
